$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new sample content. Order matters for shared-string table layout:
# G2 first, then the three header relabels, then the two long D2/F2 notes.
$ws.Range("G2").Value = "HW * Exam"
$ws.Range("D1").Value = "topic"
$ws.Range("F1").Value = "prep"
$ws.Range("G1").Value = "eval"
$ws.Range("D2").Value = "This is a test `n* This is just a test"
$ws.Range("F2").Value = "Do this * `nand this"

# Match the wrap/alignment already used by the rest of the row (style index 5).
$ws.Range("D2:G2").HorizontalAlignment = -4131
$ws.Range("D2:G2").VerticalAlignment = -4108
$ws.Range("D2:G2").WrapText = $true

# Row 2 grows to fit the new multi-line notes.
$ws.Rows("2:2").RowHeight = 31.5

# Leave the selection on F2 (single cell) instead of the whole column.
$ws.Range("F2").Select() | Out-Null
